$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Rows 53-59: tied_teams (column O) reordered
foreach ($r in 53..59) {
    $ws.Range("O$r").Value = "['Argentina', 'Colombia', 'Costa Rica', 'Ireland']"
}

# Rows 60-62: tied_teams (column O) reordered
foreach ($r in 60..62) {
    $ws.Range("O$r").Value = "['Argentina', 'Colombia']"
}

# Rows 63-73: tied_teams (column O) reordered
foreach ($r in 63..73) {
    $ws.Range("O$r").Value = "['Argentina', 'Colombia', 'Scotland', 'Austria']"
}

# Row 78: tied_teams (column O) reordered
$ws.Range("O78").Value = "['South Korea', 'Netherlands']"

# Rows 109-110: Argentina -> Bulgaria in Group D (column J) and top_four (column M)
$ws.Range("J109").Value = "['Bulgaria', 6, 3, 6]"
$ws.Range("M109").Value = "['Bulgaria', 'Belgium', 'United States', 'Italy']"
$ws.Range("P109").Value = 0
$ws.Range("Q109").Value = 11

$ws.Range("J110").Value = "['Bulgaria', 6, 3, 6]"
$ws.Range("M110").Value = "['Bulgaria', 'Belgium', 'United States', 'Italy']"
$ws.Range("Q110").Value = 11
